$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.37"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'43.95"
$ws.Range("E3").Value = "'0.54%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.578"
$ws.Range("E4").Value = "'1.17%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.08033"
$ws.Range("E5").Value = "'-0.51%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'1.976"
$ws.Range("E6").Value = "'5.11%"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'4.336"
$ws.Range("E7").Value = "'1.25%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'0.9483"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'2.548"
$ws.Range("E9").Value = "'-8.00%"
$ws.Range("G9").Value = "'22"
$ws.Range("E10").Value = "'0.47%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.1849"
$ws.Range("E11").Value = "'-2.14%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'11.83"
$ws.Range("E12").Value = "'38.44%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.09790"
$ws.Range("E13").Value = "'2.23%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.04712"
$ws.Range("E14").Value = "'14.71%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'0.16%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.001288"
$ws.Range("E16").Value = "'1.04%"
$ws.Range("G16").Value = "'22"
$ws.Range("D17").Value = "'0.04229"
$ws.Range("E17").Value = "'-2.54%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'0.005964"
$ws.Range("E18").Value = "'0.71%"
$ws.Range("G18").Value = "'22"
$ws.Range("D19").Value = "'3.370"
$ws.Range("E19").Value = "'-5.52%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.3474"
$ws.Range("E20").Value = "'-0.27%"
$ws.Range("G20").Value = "'22"
$ws.Range("D21").Value = "'0.1410"
$ws.Range("E21").Value = "'3.36%"
$ws.Range("G21").Value = "'22"
$ws.Range("E22").Value = "'-3.08%"
$ws.Range("G22").Value = "'22"
$ws.Range("D23").Value = "'0.001255"
$ws.Range("E23").Value = "'1.64%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.004304"
$ws.Range("E24").Value = "'-5.24%"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.0001193"
$ws.Range("E25").Value = "'-3.14%"
$ws.Range("G25").Value = "'22"
$ws.Range("E26").Value = "'-0.40%"
$ws.Range("G26").Value = "'22"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("D38").Value = "'0.02594"
$ws.Range("E38").Value = "'-1.65%"
$ws.Range("G38").Value = "'22"
$ws.Range("D39").Value = "'0.05511"
$ws.Range("E39").Value = "'1.57%"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.007574"
$ws.Range("E40").Value = "'-1.31%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.1401"
$ws.Range("E41").Value = "'0.74%"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.007978"
$ws.Range("E42").Value = "'-27.64%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002021"
$ws.Range("E43").Value = "'-4.26%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.008382"
$ws.Range("E44").Value = "'-9.33%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00007102"
$ws.Range("E45").Value = "'1.34%"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("E46").Value = "'0.12%"
$ws.Range("G46").Value = "'22"
$ws.Range("E47").Value = "'1.26%"
$ws.Range("G47").Value = "'22"
$ws.Range("D48").Value = "'0.004840"
$ws.Range("E48").Value = "'35.81%"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.00002106"
$ws.Range("E49").Value = "'0.12%"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.0002006"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("G50").Value = "'22"
$ws.Range("G51").Value = "'22"
